$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Usage" column (C) for Bus and Truck rows to only allow Corporate
$ws.Range("C3").Value = "[Corporate]"
$ws.Range("C4").Value = "[Corporate]"

# Update "Cov. Sum" column (E) step size from 25000 to 100000 for all data rows
$ws.Range("E2").Value = "500000-1000000,100000"
$ws.Range("E3").Value = "500000-1000000,100000"
$ws.Range("E4").Value = "500000-1000000,100000"

# Update "Value" column (F) parameters per row
$ws.Range("F2").Value = "10000-50000,10000"
$ws.Range("F3").Value = "20000-200000,40000"
$ws.Range("F4").Value = "30000-2500000,100000"

# Update selected cell to match new selection in the saved workbook
$ws.Range("C2").Select() | Out-Null
